$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3 stays a single space " " (unchanged value, already present)
$ws.Range("C3").Value = " "

# A2 becomes a text value "0002" - apply text format first so leading zero is preserved
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "0002"

# C2 becomes a text timestamp string, format stays as Text (style 1 already applied there)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2025-08-18 00:01:00"

# Update header row
$ws.Range("A1").Value = "ID_scada"
$ws.Range("B1").Value = "valor"
$ws.Range("C1").Value = "timestamp"

# Update the selected cell / active cell to C1
$ws.Range("C1").Select()
